# Edit script: Fruta / hortaliza, semanal
# Adds a week's worth of new records (5 new rows) to the "Damasco" sheet for
# "Vega Modelo de Temuco". Row 55 stays untouched; five brand-new rows are
# inserted right after it (becoming the new rows 56-60), and everything that
# used to be rows 56-67 shifts down, unchanged, to rows 61-72 - matching the
# new dimension A1:T72.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new blank rows starting at row 56. Doing this 5 times at the same
# index pushes the prior contents of row 56 (and everything below) down by
# five rows, ending up at row 61.
$ws.Rows.Item(56).Insert()
$ws.Rows.Item(56).Insert()
$ws.Rows.Item(56).Insert()
$ws.Rows.Item(56).Insert()
$ws.Rows.Item(56).Insert()

# --- Row 56: new record (Castle Brite / Especial) ---
$ws.Range("A56").Value = 10
$ws.Range("B56").Value = "Vega Modelo de Temuco"
$ws.Range("C56").Value = "La Araucanía"
$ws.Range("D56").Value = 44900
$ws.Range("E56").Value = 9
$ws.Range("F56").Value = "Fruta"
$ws.Range("G56").Value = 100103
$ws.Range("H56").Value = "Frutos de hueso (carozo)"
$ws.Range("I56").Value = 100103003
$ws.Range("J56").Value = "Damasco"
$ws.Range("K56").Value = "Castle Brite"
$ws.Range("L56").Value = "Especial"
$ws.Range("M56").Value = 50
$ws.Range("N56").Value = 27000
$ws.Range("O56").Value = 27000
$ws.Range("P56").Value = 27000
$ws.Range("Q56").Value = "$/bandeja 18 kilos"
$ws.Range("R56").Value = "Región de O'Higgins"
$ws.Range("S56").Value = 1500
$ws.Range("T56").Value = 18

# --- Row 57: new record (Castle Brite / Especial) ---
$ws.Range("A57").Value = 10
$ws.Range("B57").Value = "Vega Modelo de Temuco"
$ws.Range("C57").Value = "La Araucanía"
$ws.Range("D57").Value = 44900
$ws.Range("E57").Value = 9
$ws.Range("F57").Value = "Fruta"
$ws.Range("G57").Value = 100103
$ws.Range("H57").Value = "Frutos de hueso (carozo)"
$ws.Range("I57").Value = 100103003
$ws.Range("J57").Value = "Damasco"
$ws.Range("K57").Value = "Castle Brite"
$ws.Range("L57").Value = "Especial"
$ws.Range("M57").Value = 200
$ws.Range("N57").Value = 22000
$ws.Range("O57").Value = 22000
$ws.Range("P57").Value = 22000
$ws.Range("Q57").Value = "$/caja 18 kilos"
$ws.Range("R57").Value = "Región de O'Higgins"
$ws.Range("S57").Value = 1222
$ws.Range("T57").Value = 18

# --- Row 58: new record (Castle Brite / Primera) ---
$ws.Range("A58").Value = 10
$ws.Range("B58").Value = "Vega Modelo de Temuco"
$ws.Range("C58").Value = "La Araucanía"
$ws.Range("D58").Value = 44900
$ws.Range("E58").Value = 9
$ws.Range("F58").Value = "Fruta"
$ws.Range("G58").Value = 100103
$ws.Range("H58").Value = "Frutos de hueso (carozo)"
$ws.Range("I58").Value = 100103003
$ws.Range("J58").Value = "Damasco"
$ws.Range("K58").Value = "Castle Brite"
$ws.Range("L58").Value = "Primera"
$ws.Range("M58").Value = 480
$ws.Range("N58").Value = 20000
$ws.Range("O58").Value = 25000
$ws.Range("P58").Value = 23062
$ws.Range("Q58").Value = "$/bandeja 18 kilos"
$ws.Range("R58").Value = "Región de O'Higgins"
$ws.Range("S58").Value = 1281
$ws.Range("T58").Value = 18

# --- Row 59: new record (Castle Brite / Primera) ---
$ws.Range("A59").Value = 10
$ws.Range("B59").Value = "Vega Modelo de Temuco"
$ws.Range("C59").Value = "La Araucanía"
$ws.Range("D59").Value = 44900
$ws.Range("E59").Value = 9
$ws.Range("F59").Value = "Fruta"
$ws.Range("G59").Value = 100103
$ws.Range("H59").Value = "Frutos de hueso (carozo)"
$ws.Range("I59").Value = 100103003
$ws.Range("J59").Value = "Damasco"
$ws.Range("K59").Value = "Castle Brite"
$ws.Range("L59").Value = "Primera"
$ws.Range("M59").Value = 250
$ws.Range("N59").Value = 22000
$ws.Range("O59").Value = 22000
$ws.Range("P59").Value = 22000
$ws.Range("Q59").Value = "$/caja 18 kilos"
$ws.Range("R59").Value = "Región de O'Higgins"
$ws.Range("S59").Value = 1222
$ws.Range("T59").Value = 18

# --- Row 60: new record (Castle Brite / Segunda) ---
$ws.Range("A60").Value = 10
$ws.Range("B60").Value = "Vega Modelo de Temuco"
$ws.Range("C60").Value = "La Araucanía"
$ws.Range("D60").Value = 44900
$ws.Range("E60").Value = 9
$ws.Range("F60").Value = "Fruta"
$ws.Range("G60").Value = 100103
$ws.Range("H60").Value = "Frutos de hueso (carozo)"
$ws.Range("I60").Value = 100103003
$ws.Range("J60").Value = "Damasco"
$ws.Range("K60").Value = "Castle Brite"
$ws.Range("L60").Value = "Segunda"
$ws.Range("M60").Value = 100
$ws.Range("N60").Value = 20000
$ws.Range("O60").Value = 20000
$ws.Range("P60").Value = 20000
$ws.Range("Q60").Value = "$/bandeja 18 kilos"
$ws.Range("R60").Value = "Región de O'Higgins"
$ws.Range("S60").Value = 1111
$ws.Range("T60").Value = 18

# Rows 61-72 now hold exactly what used to be rows 56-67, unchanged, having
# simply been shifted down by five by the insert operations above, so
# nothing else needs to be written.
